$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed) date column C for all data rows (2..472)
# from 45205 (2023-10-06) to 45206 (2023-10-07)
$ws.Range("C2:C472").Value = 45206

# Row 472 gains an explicit row height (ht="15" customHeight="1") in the new file
$ws.Rows.Item(472).RowHeight = 15

# Append the new row 473 with the new cleavage notice entry
$ws.Range("A473").Value = "A 47998-2023"

$ws.Range("B473").Value = 45204
$ws.Range("B473").NumberFormat = "YYYY-MM-DD"

$ws.Range("C473").Value = 45206
$ws.Range("C473").NumberFormat = "YYYY-MM-DD"

$ws.Range("D473").Value = "HALLANDS LÄN"
$ws.Range("E473").Value = "HYLTE"

$ws.Range("G473").Value = 1
$ws.Range("H473").Value = 0
$ws.Range("I473").Value = 0
$ws.Range("J473").Value = 0
$ws.Range("K473").Value = 0
$ws.Range("L473").Value = 0
$ws.Range("M473").Value = 0
$ws.Range("N473").Value = 0
$ws.Range("O473").Value = 0
$ws.Range("P473").Value = 0
$ws.Range("Q473").Value = 0

$ws.Range("R473").Value = ""
$ws.Range("R473").WrapText = $true
